$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($Row, $Col, $Val)
    $cell = $ws.Cells.Item($Row, $Col)
    $cell.NumberFormat = "@"
    $cell.Value = $Val
    $cell.Style = "Normal"
}

Set-TextValue 2 4 "248.51"
Set-TextValue 2 7 "5"
Set-TextValue 3 4 "21.79"
Set-TextValue 3 7 "5"
Set-TextValue 4 4 "5.511"
Set-TextValue 4 7 "5"
Set-TextValue 5 4 "0.05653"
Set-TextValue 5 7 "5"
Set-TextValue 6 4 "3.379"
Set-TextValue 6 7 "5"
Set-TextValue 7 4 "6.439"
Set-TextValue 7 7 "5"
Set-TextValue 8 4 "0.8021"
Set-TextValue 8 7 "5"
Set-TextValue 9 4 "1.037"
Set-TextValue 9 7 "5"
Set-TextValue 10 2 "WazirX"
Set-TextValue 10 3 "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
Set-TextValue 10 4 "0.1425"
Set-TextValue 10 5 "9WazirXWRX"
Set-TextValue 10 7 "5"
Set-TextValue 11 2 "MandalaExchangeToken"
Set-TextValue 11 3 "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
Set-TextValue 11 4 "0.07229"
Set-TextValue 11 5 "10MandalaExchangeTokenMDX"
Set-TextValue 11 7 "5"
Set-TextValue 12 2 "LiechtensteinCryptoassetsExchange"
Set-TextValue 12 3 "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
Set-TextValue 12 4 "0.03150"
Set-TextValue 12 5 "11LiechtensteinCryptoassetsExchangeLCX"
Set-TextValue 12 7 "5"
Set-TextValue 13 2 "BitrueCoin"
Set-TextValue 13 3 "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
Set-TextValue 13 4 "0.02943"
Set-TextValue 13 5 "12BitrueCoinBTR"
Set-TextValue 13 7 "5"
Set-TextValue 14 2 "BitMartToken"
Set-TextValue 14 3 "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
Set-TextValue 14 4 "0.09287"
Set-TextValue 14 5 "13BitMartTokenBMX"
Set-TextValue 14 7 "5"
Set-TextValue 15 2 "BitForexToken"
Set-TextValue 15 3 "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
Set-TextValue 15 4 "0.001657"
Set-TextValue 15 5 "14BitForexTokenBF"
Set-TextValue 15 7 "5"
Set-TextValue 16 2 "MCDex"
Set-TextValue 16 3 "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
Set-TextValue 16 4 "3.209"
Set-TextValue 16 5 "15MCDexMCB"
Set-TextValue 16 7 "5"
Set-TextValue 17 2 "CoinExToken"
Set-TextValue 17 3 "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
Set-TextValue 17 4 "0.04730"
Set-TextValue 17 5 "16CoinExTokenCET"
Set-TextValue 17 7 "5"
Set-TextValue 18 2 "One"
Set-TextValue 18 3 "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
Set-TextValue 18 4 "0.0005826"
Set-TextValue 18 5 "17OneONE"
Set-TextValue 18 7 "5"
Set-TextValue 19 4 "0.006463"
Set-TextValue 19 7 "5"
Set-TextValue 20 4 "0.005024"
Set-TextValue 20 5 "19HotbitTokenHTBBestin24h"
Set-TextValue 20 7 "5"
Set-TextValue 21 4 "0.001051"
Set-TextValue 21 7 "5"
Set-TextValue 22 7 "5"
Set-TextValue 23 4 "0.0003203"
Set-TextValue 23 7 "5"
Set-TextValue 24 4 "4.125"
Set-TextValue 24 7 "5"
Set-TextValue 25 4 "2.109"
Set-TextValue 25 7 "5"
Set-TextValue 26 7 "5"
Set-TextValue 27 7 "5"
Set-TextValue 28 7 "5"
Set-TextValue 29 7 "5"
Set-TextValue 30 7 "5"
Set-TextValue 31 7 "5"
Set-TextValue 32 7 "5"
Set-TextValue 33 7 "5"
Set-TextValue 34 7 "5"
Set-TextValue 35 7 "5"
Set-TextValue 36 7 "5"
Set-TextValue 37 7 "5"
Set-TextValue 38 7 "5"
Set-TextValue 39 7 "5"
Set-TextValue 40 4 "0.04086"
Set-TextValue 40 7 "5"
Set-TextValue 41 4 "0.006906"
Set-TextValue 41 5 "40KickTokenKICK"
Set-TextValue 41 7 "5"
Set-TextValue 42 2 "BKEXToken"
Set-TextValue 42 3 "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
Set-TextValue 42 4 "0.1044"
Set-TextValue 42 5 "41BKEXTokenBKK"
Set-TextValue 42 7 "5"
Set-TextValue 43 2 "CEJI"
Set-TextValue 43 3 "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
Set-TextValue 43 4 "0.002972"
Set-TextValue 43 5 "42CEJICEJI"
Set-TextValue 43 7 "5"
Set-TextValue 44 4 "0.009152"
Set-TextValue 44 7 "5"
Set-TextValue 45 4 "0.00005834"
Set-TextValue 45 7 "5"
Set-TextValue 46 4 "0.00000000751"
Set-TextValue 46 7 "5"
Set-TextValue 47 4 "0.7860"
Set-TextValue 47 7 "5"
Set-TextValue 48 4 "0.01708"
Set-TextValue 48 7 "5"
Set-TextValue 49 4 "0.00002102"
Set-TextValue 49 7 "5"
Set-TextValue 50 7 "5"
Set-TextValue 51 7 "5"
